# Fruta / hortaliza, semanal
# Insert a new weekly record row at row 131 (Ají - Inferno, Terminal La Palmera
# de La Serena), shifting the existing rows 131-190 down to 132-191.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 131; rows 131..190 move down to 132..191.
$ws.Rows("131").Insert()

# Populate the new row 131 with the new record's data.
$ws.Range("A131").Value = 8
$ws.Range("B131").Value = "Terminal La Palmera de La Serena"
$ws.Range("C131").Value = "Coquimbo"
$ws.Range("D131").Value = 44603
$ws.Range("E131").Value = 4
$ws.Range("F131").Value = 100112021
$ws.Range("G131").Value = "Ají"
$ws.Range("H131").Value = "Inferno"
$ws.Range("I131").Value = "Primera"
$ws.Range("J131").Value = 600
$ws.Range("K131").Value = 9500
$ws.Range("L131").Value = 10000
$ws.Range("M131").Value = 9750
$ws.Range("N131").Value = "`$/caja 15 kilos"
$ws.Range("O131").Value = "Provincia de Limarí"
$ws.Range("P131").Value = 650
$ws.Range("Q131").Value = 15
$ws.Range("R131").Value = "Hortaliza"
